$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the timestamp in A85 (was 07:37:16, corrected to 07:00:00) ---
$ws.Range("A85").Value = 45464.2916666667

# --- Append new row 86 with the latest OHLCV data point ---

# Column A: date/time, formatted like the rest of the date column
$ws.Range("A86").Font.Name = "Calibri"
$ws.Range("A86").Font.Size = 11
$ws.Range("A86").Font.Color = $ws.Range("A85").Font.Color
$ws.Range("A86").NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Range("A86").Value = 45467.5177199074

# Column B: volume
$ws.Range("B86").Value = 750

# Columns C-F: high / low / open / close
$ws.Range("C86").Value = 2.97000002861023
$ws.Range("D86").Value = 2.97000002861023
$ws.Range("E86").Value = 2.97000002861023
$ws.Range("F86").Value = 2.97000002861023

# Column G: adj_close, stored as text (matches existing shared-string data)
$ws.Range("G86").NumberFormat = "@"
$ws.Range("G86").Value = "2.97000002861023"
$ws.Range("G86").ClearFormats()

# Column H: ticker
$ws.Range("H86").Value = "ESPE.MI"
